$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Update the Temp folder path (row 19, col B) to the new RPA10 temp path
$ws.Range("B19").Value = "C:\RPA_Repository\ExPath\RPA10_Temp"

# Add new StartDate / EndDate settings rows (A48:C48, A49:C49) and their
# descriptions (C47, C48)
$ws.Range("A48").Value = "StartDate"
$ws.Range("A49").Value = "EndDate"
$ws.Range("C47").Value = "수신한 메일에서 추출한 시작날짜"
$ws.Range("C48").Value = "수신한 메일에서 추출한 종료날짜"

# Move the active selection / view to reflect the newly added rows
$ws.Range("A50").Select()
